# Uploading new data files.
# Fill in newly-recorded days_alive values and extend the sheet with an
# (empty, date-formatted) G column matching the existing D column style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New days_alive (column C) observations -------------------------------
$ws.Range("C8").Value  = 30
$ws.Range("C9").Value  = 28
$ws.Range("C23").Value = 29
$ws.Range("C24").Value = 29
$ws.Range("C35").Value = 28
$ws.Range("C43").Value = 28
$ws.Range("C46").Value = 19

# --- Extend with an empty, date-formatted G column (rows 1-30) ------------
# Copy the existing date-formatted style from D2 (style index reused, same
# as the rest of column D) onto G1:G30 so no new style entry is created.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("G1:G30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- View state: zoom + selected cell --------------------------------------
$ws.Application.ActiveWindow.Zoom = 182
$ws.Range("E23").Select() | Out-Null
